$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 4
$ws.Range("G10").Value = 3
$ws.Range("G11").Value = 4
